# Update the timesheet to reflect the new week (Feb 8 - Feb 14, 2021)
# and move the active selection from F7 to H6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Shift the week's dates in row 5 (B5:H5) forward by 141 days, i.e.
# from the week of 44094 (2020-09-20) to the week of 44235 (2021-02-08).
$ws.Range("B5").Value = 44235
$ws.Range("C5").Value = 44236
$ws.Range("D5").Value = 44237
$ws.Range("E5").Value = 44238
$ws.Range("F5").Value = 44239
$ws.Range("G5").Value = 44240
$ws.Range("H5").Value = 44241

# Move the active cell / selection to H6 (was F7).
$ws.Range("H6").Select()
